$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.668.08'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.11%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.676.74'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.94%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.01%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.28%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3909'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.19%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3941'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.33%  '

# Row 9
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.26'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.36%  '

# Row 10
$ws.Range("B10").Value = 'BinanceUSD'
$ws.Range("C10").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.002'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.10%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.392'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.16%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08662'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.89%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.16'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.56%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.308'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.90%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.753'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.78%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001313'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.51%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.683.96'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.18%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.46%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07053'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.58%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.04%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.065'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.10%  '

# Row 22
$ws.Range("E22").Value = '  +0.40%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.20%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.669.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.13%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.349'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.75%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '23.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.42%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.728'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.60%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.51'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.96%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.748'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.15%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '147.06'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.36%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.868'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.30%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.511'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.59%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.881.90'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.51%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08381'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.98%  '

# Row 35
$ws.Range("E35").Value = '  -6.38%  '

# Row 36
$ws.Range("B36").Value = 'Algorand'
$ws.Range("C36").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2819'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.68%  '

# Row 37
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.866'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.85%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9777'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.20%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09478'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.26%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.559'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.85%  '

# Row 41
$ws.Range("E41").Value = '  -3.42%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7914'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.06%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.93%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.17%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7120'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.58%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.561'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.39%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.195'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.25%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.08648'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.01%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.003'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.26%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.321'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.71%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '137.20'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.65%  '
